$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41-55 down to 42-56.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly record.
$ws.Range("A41").Value = 3
$ws.Range("B41").Value = "Femacal de La Calera"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44900
$ws.Range("E41").Value = 5
$ws.Range("F41").Value = 300000000
$ws.Range("G41").Value = "Espárragos"
$ws.Range("H41").Value = "Verde"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 950
$ws.Range("K41").Value = 1500
$ws.Range("L41").Value = 1500
$ws.Range("M41").Value = 1500
$ws.Range("N41").Value = "$/kilo"
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 1500
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"
